$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions) - update "想去人数" (interested-count) column F
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1334
$ws1.Range("F3").Value = 1208
$ws1.Range("F4").Value = 901
$ws1.Range("F7").Value = 668
$ws1.Range("F8").Value = 111
$ws1.Range("F11").Value = 2414
$ws1.Range("F12").Value = 1599
$ws1.Range("F13").Value = 1440
$ws1.Range("F14").Value = 309
$ws1.Range("F15").Value = 241
$ws1.Range("F16").Value = 585
$ws1.Range("F17").Value = 775
$ws1.Range("F18").Value = 59
$ws1.Range("F19").Value = 302
$ws1.Range("F22").Value = 23
$ws1.Range("F24").Value = 4862
$ws1.Range("F26").Value = 404
$ws1.Range("F27").Value = 69
$ws1.Range("F28").Value = 156
$ws1.Range("F29").Value = 136
$ws1.Range("F30").Value = 219
$ws1.Range("F31").Value = 96
$ws1.Range("F32").Value = 25
$ws1.Range("F33").Value = 1031
$ws1.Range("F34").Value = 705
$ws1.Range("F35").Value = 63
$ws1.Range("F36").Value = 45
$ws1.Range("F38").Value = 385
$ws1.Range("F39").Value = 1028
$ws1.Range("F40").Value = 128
$ws1.Range("F41").Value = 103
$ws1.Range("F42").Value = 165
$ws1.Range("F43").Value = 125

# Sheet 2: 演出 (Performances) - update "想去人数" column F
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 785
$ws2.Range("F12").Value = 6

# Sheet 4: 全部类型 (All types, combined) - update "想去人数" column F
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1334
$ws4.Range("F4").Value = 785
$ws4.Range("F5").Value = 1208
$ws4.Range("F6").Value = 901
$ws4.Range("F11").Value = 668
$ws4.Range("F12").Value = 111
$ws4.Range("F17").Value = 2414
$ws4.Range("F18").Value = 1599
$ws4.Range("F19").Value = 1440
$ws4.Range("F20").Value = 309
$ws4.Range("F21").Value = 241
$ws4.Range("F22").Value = 585
$ws4.Range("F24").Value = 776
$ws4.Range("F25").Value = 59
$ws4.Range("F26").Value = 302
$ws4.Range("F28").Value = 23
$ws4.Range("F29").Value = 4862
$ws4.Range("F31").Value = 404
$ws4.Range("F32").Value = 69
$ws4.Range("F33").Value = 156
$ws4.Range("F34").Value = 136
$ws4.Range("F35").Value = 219
$ws4.Range("F36").Value = 96
$ws4.Range("F37").Value = 25
$ws4.Range("F38").Value = 1031
$ws4.Range("F39").Value = 705
$ws4.Range("F40").Value = 45
$ws4.Range("F41").Value = 385
$ws4.Range("F42").Value = 1028
$ws4.Range("F43").Value = 128
$ws4.Range("F44").Value = 103
$ws4.Range("F45").Value = 165
$ws4.Range("F46").Value = 125
$ws4.Range("F49").Value = 6
